$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert a new "2022-Q1" sheet (same column layout as the other quarterly
#    sheets), cloned from "2021-Q4", positioned immediately before "总计".
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$templateSheet = $wb.Worksheets.Item("2021-Q4")
$templateSheet.Copy($totalSheet, $null)

$ws = $wb.Worksheets.Item("2021-Q4 (2)")
$ws.Name = "2022-Q1"

# The template ("2021-Q4") only has 17 data rows (rows 2-18); the new sheet
# needs 21 data rows (rows 2-22), so extend the formatting of the last
# template row down across the new rows before filling in values.
$ws.Range("A18:H18").Copy()
$ws.Range("A19:H22").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Columns that hold numeric-looking text (fund code / scale / weight figures)
# must be forced to text so they round-trip as strings, matching the source
# data (e.g. "59.01", "0.70", "004212" keep their literal formatting).
$ws.Range("B2:B22").NumberFormat = "@"
$ws.Range("D2:G22").NumberFormat = "@"

$rows = @(
    @(0, '516970', '广发中证基建工程交易型开放式指数证券投资基金', '59.01', '99.38', '2.59', '1.5284', 9),
    @(1, '420005', '天弘周期策略混合', '5.25', '89.31', '9.47', '0.4972', 2),
    @(2, '165525', '信诚中证基建工程指数（LOF）', '17.06', '94.00', '2.44', '0.4163', 9),
    @(3, '007202', '天弘优质成长企业精选混合', '4.81', '92.52', '8.37', '0.4026', 3),
    @(4, '420001', '天弘精选混合', '7.16', '71.80', '5.40', '0.3866', 2),
    @(5, '005671', '新疆前海联合研究优选灵活配置混合A', '5.62', '77.46', '5.09', '0.2861', 3),
    @(6, '002780', '新疆前海联合泓鑫灵活配置混合A', '8.50', '75.30', '3.10', '0.2635', 9),
    @(7, '011851', '天弘先进制造混合型证券投资基金A', '2.72', '91.41', '8.76', '0.2383', 1),
    @(8, '002510', '申万菱信中证500指数增强A', '4.82', '89.83', '2.15', '0.1036', 1),
    @(9, '011852', '天弘先进制造混合型证券投资基金C', '0.70', '91.41', '8.76', '0.0613', 1),
    @(10, '004694', '天弘策略精选灵活配置混合A', '1.11', '80.93', '4.43', '0.0492', 2),
    @(11, '007043', '新疆前海联合泓鑫灵活配置混合C', '1.37', '75.30', '3.10', '0.0425', 9),
    @(12, '007795', '申万菱信中证500指数增强C', '1.12', '89.83', '2.15', '0.0241', 1),
    @(13, '006478', '长盛多因子策略优选股票', '0.51', '84.41', '4.30', '0.0219', 8),
    @(14, '010253', '兴银中证500指数增强A', '2.19', '82.47', '0.96', '0.0210', 6),
    @(15, '011205', '兴银中证500指数增强C', '1.78', '82.47', '0.96', '0.0171', 6),
    @(16, '159962', '华夏中证四川国企改革ETF', '0.49', '95.82', '3.04', '0.0149', 8),
    @(17, '005672', '新疆前海联合研究优选灵活配置混合C', '0.21', '77.46', '5.09', '0.0107', 3),
    @(18, '004748', '天弘策略精选灵活配置混合C', '0.08', '80.93', '4.43', '0.0035', 2),
    @(19, '004212', '中融量化智选混合A', '0.08', '93.46', '2.73', '0.0022', 9),
    @(20, '004783', '中融量化智选混合C', '0.01', '93.46', '2.73', '0.0003', 9)
)

foreach ($row in $rows) {
    $r = [int]$row[0] + 2
    $ws.Cells.Item($r, 1).Value = [int]$row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = [int]$row[7]
}

# ---------------------------------------------------------------------------
# 2) Prepend a "2022-Q1" summary row to the "总计" sheet, pushing the
#    existing history down by one row.
# ---------------------------------------------------------------------------
$tot = $wb.Worksheets.Item("总计")
$tot.Rows(2).Insert()
$tot.Range("B2:D2").ClearFormats()
$tot.Range("A3").Copy()
$tot.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$totRows = @(
    @(0, "2022-Q1", 21, 4.39),
    @(1, "2021-Q4", 17, 2.88),
    @(2, "2021-Q3", 12, 2.87),
    @(3, "2021-Q2", 7, 0.13),
    @(4, "2021-Q1", 18, 1.18),
    @(5, "2020-Q4", 4, 0.04)
)

foreach ($row in $totRows) {
    $r = [int]$row[0] + 2
    $tot.Cells.Item($r, 1).Value = [int]$row[0]
    $tot.Cells.Item($r, 2).Value = $row[1]
    $tot.Cells.Item($r, 3).Value = $row[2]
    $tot.Cells.Item($r, 4).Value = $row[3]
}
